# Trade #14 closed at 2026-02-17 23:53:38 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade #14 (a losing trade), which nudges the
# aggregate P&L / win-rate figures down and appends the trade row itself.

$wb = $excel.ActiveWorkbook

$wsSummary  = $wb.Worksheets.Item("Summary")
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsAll      = $wb.Worksheets.Item("All Trades")
$wsMM       = $wb.Worksheets.Item("MarketMaking")

# ---------------------------------------------------------------------------
# Summary sheet - updated aggregate metrics
# ---------------------------------------------------------------------------
$wsSummary.Range("B3").Value = 1500.29   # Current Capital
$wsSummary.Range("B4").Value = 0.29      # Total P&L $
$wsSummary.Range("B5").Value = 0.41      # Total P&L %
$wsSummary.Range("B6").Value = 14        # Total Trades
$wsSummary.Range("B8").Value = 5         # Losing Trades
$wsSummary.Range("B9").Value = 57.14     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$wsStrategy.Range("C6").Value = 100.29   # Capital
$wsStrategy.Range("D6").Value = 14       # Trades
$wsStrategy.Range("E6").Value = 0.29     # P&L $
$wsStrategy.Range("F6").Value = 0.29     # P&L %
$wsStrategy.Range("G6").Value = 57.14    # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the closed trade #14 row to a trades-log sheet at row 15
# ---------------------------------------------------------------------------
function Add-Trade14Row($ws) {
    $ws.Cells.Item(15, 1).Value = 14                 # Trade #
    $ws.Cells.Item(15, 2).Value = "'2026-02-17"      # Date (force text)
    $ws.Cells.Item(15, 2).Style = "Normal"
    $ws.Cells.Item(15, 3).Value = "'23:53:32"        # Time (force text)
    $ws.Cells.Item(15, 3).Style = "Normal"
    $ws.Cells.Item(15, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item(15, 5).Value = "UP"                # Side
    $ws.Cells.Item(15, 6).Value = 0.94                # Entry Price
    $ws.Cells.Item(15, 7).Value = 0.91                # Exit Price
    $ws.Cells.Item(15, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item(15, 9).Value = -3.1915             # P&L %
    $ws.Cells.Item(15, 10).Value = -0.03              # P&L $
    $ws.Cells.Item(15, 11).Value = 100.29             # Capital After
    $ws.Cells.Item(15, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(15, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(15, 14).Value = 0.6                # Confidence
    $ws.Cells.Item(15, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(15, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item(15, 17).Value = 0.15               # Duration (min)
}

# ---------------------------------------------------------------------------
# All Trades sheet - append trade #14
# ---------------------------------------------------------------------------
Add-Trade14Row $wsAll

# ---------------------------------------------------------------------------
# MarketMaking sheet - append trade #14 (mirrors All Trades)
# ---------------------------------------------------------------------------
Add-Trade14Row $wsMM
